$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows after row 16 (pushing the old row 17 "sum" row down to row 20).
# Copy formatting down from row 16 so the new rows inherit the same look.
$ws.Rows("17:19").Insert()

# Match formatting of the template row (16) for the three new data rows
# (values are filled in afterwards so this paste only carries formats).
$ws.Range("A16:H16").Copy()
$ws.Range("A17:H19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 17: Front Page
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Front Page"
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 0.27083333333333331
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0

# Row 18: Forum
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Forum"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 0.40625
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0

# Row 19: Profile
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Profile"
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 0.1875
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0

# Merge G:H on each new row, same as the existing rows above.
$ws.Range("G17:H17").Merge()
$ws.Range("G18:H18").Merge()
$ws.Range("G19:H19").Merge()

# Row 20 (previously row 17) totals now span the four data rows.
$ws.Range("C20").Formula = "=SUM(C16,C17,C18,C19)"
$ws.Range("D20").Formula = "=SUM(D16,D17,D18,D19)"
$ws.Range("E20").Formula = "=COUNT(E16:E20)"
$ws.Range("F20").Formula = "=COUNT(F16:F20)"

$ws.Calculate()

$ws.Range("F25").Select()
